$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9108380675315857
$ws.Range("B1").Value = 1.643423914909363
$ws.Range("C1").Value = 4.347991466522217
$ws.Range("D1").Value = 2.603734731674194
$ws.Range("E1").Value = 0.803712010383606
